# Rule factoring.xlsx update (#1936)
# Applies the content changes described by the commit's unified diff:
#  - F5:  IN PROGRESS -> TODO
#  - F6:  IN PROGRESS -> DONE
#  - F22: TODO        -> DONE
#  - F37: IN PROGRESS -> DONE
#  - F38: IN PROGRESS -> DONE
#  - F49: TODO        -> IN PROGRESS: EN
#  - F62: TODO        -> IN PROGRESS: EN
#  - F63: TODO        -> IN PROGRESS: EN
#  - H35: new note "Talk to MF. Eddy already has a PR for 1022 with more checks."
#  - I2 (COUNTIF formula) recalculates to the new count (12)
#  - Window zoomed to 190% with E47 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New note cell first, so its shared string gets the lower of the two
#     new shared-string indices (matches the order they appear in the file) ---
$ws.Range("H35").Value = "Talk to MF. Eddy already has a PR for 1022 with more checks."
$ws.Range("H35").Font.Bold = $true

# --- Update status cells in column F, reusing the same cell formatting that
#     is already used elsewhere in the sheet for the same status value ---

# TODO = red fill (style already used by F8)
$ws.Range("F8").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# DONE = green fill (style already used by F4)
$ws.Range("F4").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F38").PasteSpecial(-4122)

# IN PROGRESS: EN = theme fill (style already used by G22's "IN PROGRESS")
$ws.Range("G22").Copy()
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("G22").Copy()
$ws.Range("F62").PasteSpecial(-4122)
$ws.Range("G22").Copy()
$ws.Range("F63").PasteSpecial(-4122)

$ws.Range("F5").Value = "TODO"
$ws.Range("F6").Value = "DONE"
$ws.Range("F22").Value = "DONE"
$ws.Range("F37").Value = "DONE"
$ws.Range("F38").Value = "DONE"
$ws.Range("F49").Value = "IN PROGRESS: EN"
$ws.Range("F62").Value = "IN PROGRESS: EN"
$ws.Range("F63").Value = "IN PROGRESS: EN"

# --- Recalculate the TODO/IN PROGRESS counter now that the source cells changed ---
$ws.Range("I2").Formula = '=COUNTIF(F:F, "TODO") + COUNTIF(F:F, "IN PROGRESS")'

# --- Window view: zoom + selection (matches the author's last saved view) ---
$win = $excel.ActiveWindow
$win.Zoom = 190
$ws.Range("E47").Select()
